$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    ,@(17, 'Luyện tập cấu trúc điều kiện 01 - Bài 1', 'https://github.com/thanhviet05x1d/c0523g1_nguyen_thanh_viet_module1/blob/main/s09_dieu_kien_2/%5BB%C3%A0i%20t%E1%BA%ADp%209-1%5D%20Luy%E1%BB%87n%20t%E1%BA%ADp%20c%E1%BA%A5u%20tr%C3%BAc%20%C4%91i%E1%BB%81u%20ki%E1%BB%87n%2001.html')
    ,@(18, 'Luyện tập cấu trúc điều kiện 01 - Bài 2', 'https://github.com/thanhviet05x1d/c0523g1_nguyen_thanh_viet_module1/blob/main/s09_dieu_kien_2/%5BB%C3%A0i%20t%E1%BA%ADp%209-2%5D%20Luy%E1%BB%87n%20t%E1%BA%ADp%20c%E1%BA%A5u%20tr%C3%BAc%20%C4%91i%E1%BB%81u%20ki%E1%BB%87n%2001.html')
    ,@(19, 'Luyện tập cấu trúc điều kiện 01 - Bài 3', 'https://github.com/thanhviet05x1d/c0523g1_nguyen_thanh_viet_module1/blob/main/s09_dieu_kien_2/%5BB%C3%A0i%20t%E1%BA%ADp%209-3%5D%20Luy%E1%BB%87n%20t%E1%BA%ADp%20c%E1%BA%A5u%20tr%C3%BAc%20%C4%91i%E1%BB%81u%20ki%E1%BB%87n%2001.html')
    ,@(20, 'Luyện tập cấu trúc điều kiện 01 - Bài 4', 'https://github.com/thanhviet05x1d/c0523g1_nguyen_thanh_viet_module1/blob/main/s09_dieu_kien_2/%5BB%C3%A0i%20t%E1%BA%ADp%209-4%5D%20Luy%E1%BB%87n%20t%E1%BA%ADp%20c%E1%BA%A5u%20tr%C3%BAc%20%C4%91i%E1%BB%81u%20ki%E1%BB%87n%2001.html')
    ,@(21, 'Luyện tập cấu trúc điều kiện 01 - Bài 5', 'https://github.com/thanhviet05x1d/c0523g1_nguyen_thanh_viet_module1/blob/main/s09_dieu_kien_2/%5BB%C3%A0i%20t%E1%BA%ADp%209-5%5D%20Luy%E1%BB%87n%20t%E1%BA%ADp%20c%E1%BA%A5u%20tr%C3%BAc%20%C4%91i%E1%BB%81u%20ki%E1%BB%87n%2001.html')
    ,@(22, 'Luyện tập cấu trúc điều kiện 01 - Bài 6', 'https://github.com/thanhviet05x1d/c0523g1_nguyen_thanh_viet_module1/blob/main/s09_dieu_kien_2/%5BB%C3%A0i%20t%E1%BA%ADp%209-6%5D%20Luy%E1%BB%87n%20t%E1%BA%ADp%20c%E1%BA%A5u%20tr%C3%BAc%20%C4%91i%E1%BB%81u%20ki%E1%BB%87n%2001.html')
    ,@(23, 'Luyện tập cấu trúc điều kiện 01 - Bài 7', 'https://github.com/thanhviet05x1d/c0523g1_nguyen_thanh_viet_module1/blob/main/s09_dieu_kien_2/%5BB%C3%A0i%20t%E1%BA%ADp%209-7%5D%20Luy%E1%BB%87n%20t%E1%BA%ADp%20c%E1%BA%A5u%20tr%C3%BAc%20%C4%91i%E1%BB%81u%20ki%E1%BB%87n%2001.html')
    ,@(24, 'Luyện tập cấu trúc điều kiện 02 - Bài 1', 'https://github.com/thanhviet05x1d/c0523g1_nguyen_thanh_viet_module1/blob/main/s09_dieu_kien_2/%5BB%C3%A0i%20t%E1%BA%ADp%209-8%5D%20Luy%E1%BB%87n%20t%E1%BA%ADp%20c%E1%BA%A5u%20tr%C3%BAc%20%C4%91i%E1%BB%81u%20ki%E1%BB%87n%2002.html')
    ,@(25, 'Luyện tập cấu trúc điều kiện 02 - Bài 2', 'https://github.com/thanhviet05x1d/c0523g1_nguyen_thanh_viet_module1/blob/main/s09_dieu_kien_2/%5BB%C3%A0i%20t%E1%BA%ADp%209-9%5D%20Luy%E1%BB%87n%20t%E1%BA%ADp%20c%E1%BA%A5u%20tr%C3%BAc%20%C4%91i%E1%BB%81u%20ki%E1%BB%87n%2002.html')
    ,@(26, 'Luyện tập cấu trúc điều kiện 02 - Bài 3', 'https://github.com/thanhviet05x1d/c0523g1_nguyen_thanh_viet_module1/blob/main/s09_dieu_kien_2/%5BB%C3%A0i%20t%E1%BA%ADp%209-10%5D%20Luy%E1%BB%87n%20t%E1%BA%ADp%20c%E1%BA%A5u%20tr%C3%BAc%20%C4%91i%E1%BB%81u%20ki%E1%BB%87n%2002.html')
    ,@(27, 'Luyện tập cấu trúc điều kiện 02 - Bài 4', 'https://github.com/thanhviet05x1d/c0523g1_nguyen_thanh_viet_module1/blob/main/s09_dieu_kien_2/%5BB%C3%A0i%20t%E1%BA%ADp%209-11%5D%20Luy%E1%BB%87n%20t%E1%BA%ADp%20c%E1%BA%A5u%20tr%C3%BAc%20%C4%91i%E1%BB%81u%20ki%E1%BB%87n%2002.html')
    ,@(28, 'Luyện tập cấu trúc điều kiện 02 - Bài 5', 'https://github.com/thanhviet05x1d/c0523g1_nguyen_thanh_viet_module1/blob/main/s09_dieu_kien_2/%5BB%C3%A0i%20t%E1%BA%ADp%209-12%5D%20Luy%E1%BB%87n%20t%E1%BA%ADp%20c%E1%BA%A5u%20tr%C3%BAc%20%C4%91i%E1%BB%81u%20ki%E1%BB%87n%2002.html')
    ,@(29, 'Luyện tập cấu trúc điều kiện 02 - Bài 6', 'https://github.com/thanhviet05x1d/c0523g1_nguyen_thanh_viet_module1/blob/main/s09_dieu_kien_2/%5BB%C3%A0i%20t%E1%BA%ADp%209-13%5D%20Luy%E1%BB%87n%20t%E1%BA%ADp%20c%E1%BA%A5u%20tr%C3%BAc%20%C4%91i%E1%BB%81u%20ki%E1%BB%87n%2002.html')
    ,@(30, 'Luyện tập cấu trúc điều kiện 02 - Bài 7', 'https://github.com/thanhviet05x1d/c0523g1_nguyen_thanh_viet_module1/blob/main/s09_dieu_kien_2/%5BB%C3%A0i%20t%E1%BA%ADp%209-14%5D%20Luy%E1%BB%87n%20t%E1%BA%ADp%20c%E1%BA%A5u%20tr%C3%BAc%20%C4%91i%E1%BB%81u%20ki%E1%BB%87n%2002.html')
    ,@(31, 'Luyện tập cấu trúc điều kiện 02 - Bài 8', 'https://github.com/thanhviet05x1d/c0523g1_nguyen_thanh_viet_module1/blob/main/s09_dieu_kien_2/%5BB%C3%A0i%20t%E1%BA%ADp%209-15%5D%20Luy%E1%BB%87n%20t%E1%BA%ADp%20c%E1%BA%A5u%20tr%C3%BAc%20%C4%91i%E1%BB%81u%20ki%E1%BB%87n%2002.html')
    ,@(32, 'Luyện tập cấu trúc điều kiện 02 - Bài 9', 'https://github.com/thanhviet05x1d/c0523g1_nguyen_thanh_viet_module1/blob/main/s09_dieu_kien_2/%5BB%C3%A0i%20t%E1%BA%ADp%209-16%5D%20Luy%E1%BB%87n%20t%E1%BA%ADp%20c%E1%BA%A5u%20tr%C3%BAc%20%C4%91i%E1%BB%81u%20ki%E1%BB%87n%2002.html')
    ,@(33, 'Luyện tập cấu trúc điều kiện 02 - Bài 10', 'https://github.com/thanhviet05x1d/c0523g1_nguyen_thanh_viet_module1/blob/main/s09_dieu_kien_2/%5BB%C3%A0i%20t%E1%BA%ADp%209-17%5D%20Luy%E1%BB%87n%20t%E1%BA%ADp%20c%E1%BA%A5u%20tr%C3%BAc%20%C4%91i%E1%BB%81u%20ki%E1%BB%87n%2002.html')
    ,@(34, 'Luyện tập cấu trúc điều kiện 02 - Bài 11', 'https://github.com/thanhviet05x1d/c0523g1_nguyen_thanh_viet_module1/blob/main/s09_dieu_kien_2/%5BB%C3%A0i%20t%E1%BA%ADp%209-18%5D%20Luy%E1%BB%87n%20t%E1%BA%ADp%20c%E1%BA%A5u%20tr%C3%BAc%20%C4%91i%E1%BB%81u%20ki%E1%BB%87n%2002.html')
    ,@(35, 'Luyện tập cấu trúc điều kiện 02 - Bài 12', 'https://github.com/thanhviet05x1d/c0523g1_nguyen_thanh_viet_module1/blob/main/s09_dieu_kien_2/%5BB%C3%A0i%20t%E1%BA%ADp%209-19%5D%20Luy%E1%BB%87n%20t%E1%BA%ADp%20c%E1%BA%A5u%20tr%C3%BAc%20%C4%91i%E1%BB%81u%20ki%E1%BB%87n%2002.html')
)

foreach ($item in $data) {
    $row = $item[0]
    $ws.Cells.Item($row, 1).Value = $item[1]
    $ws.Cells.Item($row, 2).Value = $item[2]
}

$ws.Range("B17:B35").Select()
$excel.ActiveWindow.ScrollRow = 13
